$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.866.26'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").Value = '2.320.71'
$ws.Range("E3").Value = '  +1.75%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '114.20'
$ws.Range("E5").Value = '  +19.24%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '270.58'
$ws.Range("E6").Value = '  +1.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.623'
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("E8").Value = '  +0.21%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.626'
$ws.Range("E9").Value = '  +2.87%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '47.48'
$ws.Range("E10").Value = '  +7.36%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0946'
$ws.Range("E11").Value = '  +1.14%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.84'
$ws.Range("E12").Value = '  +14.41%  '

$ws.Range("E13").Value = '  +2.12%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.81'
$ws.Range("E14").Value = '  +4.32%  '

$ws.Range("D15").Value = '2.664.50'
$ws.Range("E15").Value = '  +2.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.865'
$ws.Range("E16").Value = '  +1.92%  '

$ws.Range("D17").Value = '2.314.48'
$ws.Range("E17").Value = '  +1.31%  '

$ws.Range("D18").Value = '43.885.64'
$ws.Range("E18").Value = '  +0.73%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000111'
$ws.Range("E19").Value = '  +3.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.71'
$ws.Range("E20").Value = '  +9.21%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.69'
$ws.Range("E21").Value = '  +1.03%  '

$ws.Range("E22").Value = '  +6.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.68'
$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.88'
$ws.Range("E24").Value = '  +15.75%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.47'
$ws.Range("E25").Value = '  +5.85%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.50'
$ws.Range("E27").Value = '  +1.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '42.70'
$ws.Range("E28").Value = '  +11.38%  '

$ws.Range("E29").Value = '  -0.13%  '

$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '178.03'
$ws.Range("E31").Value = '  +1.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.01'
$ws.Range("E32").Value = '  +0.03%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0937'
$ws.Range("E33").Value = '  +6.08%  '

$ws.Range("E34").Value = '  +4.37%  '

$ws.Range("B35").Value = 'RenderToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.75'
$ws.Range("E35").Value = '  +6.96%  '

$ws.Range("B36").Value = 'Stellar'
$ws.Range("C36").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.125'
$ws.Range("E36").Value = '  -0.30%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.112'
$ws.Range("E37").Value = '  +3.94%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0358'
$ws.Range("E38").Value = '  +0.81%  '

$ws.Range("B39").Value = 'NEARProtocol'
$ws.Range("C39").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.89'
$ws.Range("E39").Value = '  +18.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.248'
$ws.Range("E40").Value = '  +4.13%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.41'
$ws.Range("E41").Value = '  +2.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '69.92'
$ws.Range("E42").Value = '  +11.13%  '

$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.41'
$ws.Range("E43").Value = '  +4.59%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.08%  '

$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.65'
$ws.Range("E45").Value = '  +6.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.92'
$ws.Range("E46").Value = '  +13.49%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.82'
$ws.Range("E47").Value = '  -0.15%  '

$ws.Range("E48").Value = '  -0.76%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '100.44'
$ws.Range("E49").Value = '  +2.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.466'
$ws.Range("E50").Value = '  +10.26%  '

$ws.Range("E51").Value = '  +3.17%  '
